$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values include strings that look numeric (e.g. "5.960", "0.00001046")
# or that use multiple dots as thousands separators (e.g. "27.697.25"). Force
# Text formatting before assignment so Excel stores the exact literal digits/
# trailing zeros instead of re-parsing them as a number (which would drop
# trailing zeros or switch to scientific notation), then restore the default
# "Normal" style so the cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.697.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4721"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.99%  "
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08051"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.026"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.960"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.147"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001046"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06648"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.709.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.522"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.100.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.099"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.603"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9738"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09551"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.593"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.340"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02256"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.231"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.247"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6029"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1895"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.256"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5696"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.950"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.375"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06872"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("E51").Value = "  +14.51%  "
